$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.324.67'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.869.46'
$ws.Range("E3").Value = '  +0.50%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.13%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.7102'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.81%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '241.07'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.20%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.14%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.07871'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.57%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.3079'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.06%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '25.37'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +3.20%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.08248'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.871.15'
$ws.Range("E12").Value = '  -25.33%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.7215'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.234'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.17%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '90.65'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").Value = '29.318.19'
$ws.Range("E16").Value = '  +1.46%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '5.839'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.39%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '244.12'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.38%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000007815'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.21%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.19'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").Value = '2.116.28'
$ws.Range("E21").Value = '  +8.81%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '8.007'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +6.53%  '
$ws.Range("E24").Value = '  +0.15%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.1591'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +11.09%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '162.43'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.14%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.983'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.95%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '18.22'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.64%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.354'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.81%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.495'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.43%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.379'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.02%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.087'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.67%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.05180'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +0.68%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.186'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.75%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.7199'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.20%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.671'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.11%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01854'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.17%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.690'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Value = '1.174.03'
$ws.Range("E40").Value = '  +0.72%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.9027'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.41%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '6.104'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +2.66%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '72.49'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.34%  '
$ws.Range("E44").Value = '  +0.22%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '101.92'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.63%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.5284'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.79%  '
$ws.Range("D47").Value = '2.008.09'
$ws.Range("E47").Value = '  +8.43%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.787'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.66%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.897'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +5.78%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '9.250'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.00000000118'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.31%  '
